$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing values in rows 1-6 (B and E columns; formulas recalc automatically) ----
$ws.Range("B1").Value = 673
$ws.Range("E1").Value = 1205

$ws.Range("B2").Value = 83

$ws.Range("B3").Value = 48

$ws.Range("B4").Value = 35
$ws.Range("E4").Value = 81

$ws.Range("B5").Value = 64
$ws.Range("E5").Value = 3

$ws.Range("B6").Value = 15

# ---- New block: rows 12-18, usage hours table ----
$ws.Range("A12").Value = "18-19"
$ws.Range("B12").Value = 10.5
$ws.Range("C12").Value = 10.5
$ws.Range("D12").Value = 2
$ws.Range("E12").Value = 0.5

$ws.Range("A13").Value = "19-20"
$ws.Range("B13").Value = 12.2
$ws.Range("C13").Value = 12.2
$ws.Range("D13").Value = 2.4
$ws.Range("E13").Value = 0.6

$ws.Range("A14").Value = "20-21"
$ws.Range("B14").Value = 14.5
$ws.Range("C14").Value = 14.5
$ws.Range("D14").Value = 2.9
$ws.Range("E14").Value = 0.9

$ws.Range("A15").Value = "21-22"
$ws.Range("B15").Value = 19
$ws.Range("C15").Value = 19
$ws.Range("D15").Value = 3.4
$ws.Range("E15").Value = 1.1

$ws.Range("A16").Value = "22-23"
$ws.Range("B16").Value = 25
$ws.Range("C16").Value = 25
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 1.5

$ws.Range("A17").Value = "23-24"
$ws.Range("B17").Value = 32
$ws.Range("C17").Value = 32
$ws.Range("D17").Value = 6.4
$ws.Range("E17").Value = 2.3

$ws.Range("A18").Value = "24-25"
$ws.Range("B18").Value = 42.5
$ws.Range("C18").Value = 42.5
$ws.Range("D18").Value = 8.5
$ws.Range("E18").Value = 3

# ---- Row 21: label only (entered before row 19's label so shared-string order matches) ----
$ws.Range("A21").Value = "현재"

# ---- Totals row 19 ----
$ws.Range("A19").Value = "성성"
$ws.Range("B19").Formula = "=SUM(B12:B18)"
$ws.Range("C19").Formula = "=SUM(C12:C18)*2"
$ws.Range("D19").Formula = "=SUM(D12:D18)*2"
$ws.Range("E19").Formula = "=SUM(E12:E18)*2"

# ---- Row 20: reference capacity numbers ----
$ws.Range("B20").Value = 200
$ws.Range("C20").Value = 350
$ws.Range("D20").Value = 70
$ws.Range("E20").Value = 30

# ---- Row 22: remaining capacity ----
$ws.Range("B22").Formula = "=B20-B21"
$ws.Range("C22:E22").Formula = "=C20-C21"

# ---- Row 23: per-15-minute value, formatted with 2 decimals ----
$ws.Range("B23:E23").NumberFormat = "0.00"
$ws.Range("B23").Formula = "=B22/15"
$ws.Range("C23:E23").Formula = "=C22/15"

# ---- Row 24: empty cells carrying the same number format ----
$ws.Range("B24:E24").NumberFormat = "0.00"

# ---- Update sheet view to match final selection/scroll position ----
$ws.Application.ActiveWindow.ScrollRow = 14
[void]$ws.Range("G22").Select()
